$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spezpreise")

# Copy the formatting of the existing data row down into the new row so the
# new cells reuse the same styles (date format / currency format) instead of
# creating brand-new number formats.
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("D2").Copy($ws.Range("D3"))

# Fill in the new row's values.
$ws.Range("A3").Value = 45683
$ws.Range("B3").Value = "Spez 1"
$ws.Range("C3").Value = "Kaffee & Gipfeli"
$ws.Range("D3").Value = 5

# Keep the Excel table (ListObject) in sync with the new row.
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:D3"))

$ws.Range("C4").Select()
